$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "67.303.22"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.806.17"
$ws.Range("E3").Value = "  +6.90%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.34%  "

# Row 5
Set-TextValue $ws.Range("D5") "141.14"
$ws.Range("E5").Value = "  +6.88%  "

# Row 6
Set-TextValue $ws.Range("D6") "418.15"
$ws.Range("E6").Value = "  -0.32%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.793.39"
$ws.Range("E7").Value = "  +6.83%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.644"
$ws.Range("E8").Value = "  -2.74%  "

# Row 9
$ws.Range("E9").Value = "  +0.08%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.762"
$ws.Range("E10").Value = "  -3.25%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.177"
$ws.Range("E11").Value = "  +4.55%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000365"
$ws.Range("E12").Value = "  +25.56%  "

# Row 13
Set-TextValue $ws.Range("D13") "44.68"
$ws.Range("E13").Value = "  +2.68%  "

# Row 14
Set-TextValue $ws.Range("D14") "10.34"
$ws.Range("E14").Value = "  +1.85%  "

# Row 15
Set-TextValue $ws.Range("D15") "4.410.75"
$ws.Range("E15").Value = "  +6.91%  "

# Row 16
$ws.Range("E16").Value = "  -0.61%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.811.39"
$ws.Range("E17").Value = "  +6.35%  "

# Row 18
Set-TextValue $ws.Range("D18") "20.70"
$ws.Range("E18").Value = "  +0.97%  "

# Row 19
Set-TextValue $ws.Range("D19") "13.53"
$ws.Range("E19").Value = "  +6.27%  "

# Row 20
$ws.Range("E20").Value = "  +0.21%  "

# Row 21
Set-TextValue $ws.Range("D21") "67.469.85"
$ws.Range("E21").Value = "  +1.78%  "

# Row 22
Set-TextValue $ws.Range("D22") "437.50"
$ws.Range("E22").Value = "  -2.87%  "

# Row 23
Set-TextValue $ws.Range("D23") "15.30"
$ws.Range("E23").Value = "  +15.78%  "

# Row 24
Set-TextValue $ws.Range("D24") "89.28"
$ws.Range("E24").Value = "  -1.16%  "

# Row 25
Set-TextValue $ws.Range("D25") "3.11"
$ws.Range("E25").Value = "  -4.25%  "

# Row 26
Set-TextValue $ws.Range("D26") "37.63"
$ws.Range("E26").Value = "  +9.45%  "

# Row 27
Set-TextValue $ws.Range("D27") "3.31"
$ws.Range("E27").Value = "  -2.26%  "

# Row 28
Set-TextValue $ws.Range("D28") "9.80"
$ws.Range("E28").Value = "  -2.75%  "

# Row 29
Set-TextValue $ws.Range("D29") "5.16"
$ws.Range("E29").Value = "  +6.70%  "

# Row 30
Set-TextValue $ws.Range("D30") "12.84"
$ws.Range("E30").Value = "  +2.50%  "

# Row 31
$ws.Range("E31").Value = "  +3.65%  "

# Row 32
$ws.Range("E32").Value = "  -1.92%  "

# Row 33
Set-TextValue $ws.Range("D33") "7.19"
$ws.Range("E33").Value = "  -1.71%  "

# Row 34
Set-TextValue $ws.Range("D34") "42.28"
$ws.Range("E34").Value = "  +8.16%  "

# Row 35
$ws.Range("E35").Value = "  -0.39%  "

# Row 36
Set-TextValue $ws.Range("D36") "57.45"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
$ws.Range("E38").Value = "  -4.04%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.99"
$ws.Range("E39").Value = "  +27.63%  "

# Row 40
$ws.Range("E40").Value = "  -2.80%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0₃0681"
$ws.Range("E41").Value = "  -11.24%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.997"
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("E43").Value = "  +4.63%  "

# Row 44 (swap with row 45, plus updated values)
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D44") "2.14"
$ws.Range("E44").Value = "  +6.21%  "

# Row 45 (swap with row 44, plus updated values)
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D45") "3.20"
$ws.Range("E45").Value = "  +24.91%  "

# Row 46
Set-TextValue $ws.Range("D46") "147.48"
$ws.Range("E46").Value = "  -0.30%  "

# Row 47
$ws.Range("E47").Value = "  -1.28%  "

# Row 48
$ws.Range("E48").Value = "  -6.54%  "

# Row 49
Set-TextValue $ws.Range("D49") "26.26"
$ws.Range("E49").Value = "  +19.14%  "

# Row 50
$ws.Range("E50").Value = "  -8.23%  "

# Row 51
$ws.Range("E51").Value = "  -2.98%  "

